$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Product Data@0x4"
$ws2 = $wb.Worksheets.Item(2)   # "Statistics@0x5"

# ---------------------------------------------------------------------------
# Sheet 1: "Product Data"
# ---------------------------------------------------------------------------

# Numbers that were previously stored as text get re-entered as real numbers
$ws1.Range("E2").Value = 0
$ws1.Range("E3").Value = 0

# Hardware revision values corrected
$ws1.Range("E4").Value = 1
$ws1.Range("E5").Value = 3

# Relabel "Build" to "Revision" and fix its value
$ws1.Range("A6").Value = "Hardware Revision - Revision"
$ws1.Range("E6").Value = 1

$ws1.Range("E7").Value = 0
$ws1.Range("E8").Value = 2
$ws1.Range("E9").Value = 1
$ws1.Range("E10").Value = 4

# Release name
$ws1.Range("E11").Value = "Valerie"

# Serial number + description
$ws1.Range("E12").Value = "MyToolItStu001-1-00001-001-2"
$ws1.Range("H12").Value = "Manufactor Serial Number (Derived from ISBN); Product Group - Subgroup - Manufacture ID - Product Number - Check Digit"

# Manufacturer name
$ws1.Range("E13").Value = "Digital and Analog Communication with 5V-Supply - www.ico-tronic.com"

# OEM Free Use -> plain 0
$ws1.Range("E14").Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: "Statistics"
# ---------------------------------------------------------------------------

$ws2.Range("E2").Value = 0
$ws2.Range("E3").Value = 0
$ws2.Range("E4").Value = 0
$ws2.Range("E5").Value = 0
$ws2.Range("E6").Value = 0

# Production date corrected
$ws2.Range("E7").Value = 20190910

# Batch number row reuses the regular data-row style (instead of its own
# duplicate font/style) and its value is reset to 0
$ws2.Range("A7:H7").Copy()
$ws2.Range("A8:H8").PasteSpecial(-4122)
$ws2.Range("E8").Value = 0
